# Add a new "Week 16" blog update at the end of the document, matching
# the existing document's pattern of: two blank paragraphs, a Heading2
# paragraph with the week title, and a body paragraph with the update text.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$week16Title = "Week 16"
$week16Body  = "This week, I met with my supervisor to present my project update. I informed him about issues with the ESP32 WIFI + Camera microcontroller. He tried to help me in fixing the issue but to no avail. We later decided to use the Arduino Uno to program the ESP32 instead of the ESP32 MB programmer, and the camera module worked. I have started using the ESP 32 camera module to take pictures and test the images with some AI waste classifier models I found online."

$newContentXml = @"
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>$week16Title</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>$week16Body</w:t></w:r></w:p>
"@

# Insert at the very end of the document's main story so the existing
# final paragraph (and its text) is preserved untouched.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($newContentXml) | Out-Null
